$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 21179
$ws.Range("E2").Value = 156
$ws.Range("F2").Value = 156
$ws.Range("G2").Value = -1442
$ws.Range("H2").Value = -2211
$ws.Range("I2").Value = -2171
$ws.Range("J2").Value = -40
$ws.Range("K2").Value = 18735
$ws.Range("L2").Value = 18001
$ws.Range("M2").Value = 734
$ws.Range("N2").Value = 124
$ws.Range("O2").Value = 610
$ws.Range("P2").Value = 5196
$ws.Range("Q2").Value = -1059
$ws.Range("R2").Value = 322
$ws.Range("S2").Value = 685
$ws.Range("T2").Value = 129
$ws.Range("U2").Value = -1188
$ws.Range("V2").Value = 11373
$ws.Range("W2").Value = 0.74
$ws.Range("X2").Value = -10.44
$ws.Range("Y2").Value = -159.57
$ws.Range("Z2").Value = -10.71
$ws.Range("AA2").Value = 2451.76
$ws.Range("AB2").Value = -31.97
$ws.Range("AC2").Value = -5287
$ws.Range("AD2").Value = -1.13
$ws.Range("AE2").Value = 299
$ws.Range("AF2").Value = 20.04
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 32707337

# Row 3
$ws.Range("D3").Value = 16887
$ws.Range("E3").Value = 281
$ws.Range("F3").Value = 281
$ws.Range("G3").Value = -639
$ws.Range("H3").Value = -683
$ws.Range("I3").Value = -447
$ws.Range("J3").Value = -237
$ws.Range("K3").Value = 15528
$ws.Range("L3").Value = 11986
$ws.Range("M3").Value = 3542
$ws.Range("N3").Value = 3084
$ws.Range("O3").Value = 459
$ws.Range("P3").Value = 4195
$ws.Range("Q3").Value = -1208
$ws.Range("R3").Value = -1120
$ws.Range("S3").Value = 2821
$ws.Range("T3").Value = 82
$ws.Range("U3").Value = -1289
$ws.Range("V3").Value = 7664
$ws.Range("W3").Value = 1.66
$ws.Range("X3").Value = -4.05
$ws.Range("Y3").Value = -27.86
$ws.Range("Z3").Value = -3.99
$ws.Range("AA3").Value = 338.34
$ws.Range("AB3").Value = 56.94
$ws.Range("AC3").Value = -175
$ws.Range("AD3").Value = -13.1
$ws.Range("AE3").Value = 368
$ws.Range("AF3").Value = 6.23
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 839089231

# Row 4
$ws.Range("D4").Value = 13740
$ws.Range("E4").Value = 487
$ws.Range("F4").Value = 487
$ws.Range("G4").Value = -124
$ws.Range("H4").Value = -159
$ws.Range("I4").Value = -180
$ws.Range("J4").Value = 21
$ws.Range("K4").Value = 13642
$ws.Range("L4").Value = 10167
$ws.Range("M4").Value = 3476
$ws.Range("N4").Value = 3639
$ws.Range("O4").Value = -164
$ws.Range("P4").Value = 4282
$ws.Range("Q4").Value = -56
$ws.Range("R4").Value = 1191
$ws.Range("S4").Value = -259
$ws.Range("T4").Value = 117
$ws.Range("U4").Value = -173
$ws.Range("V4").Value = 6748
$ws.Range("W4").Value = 3.54
$ws.Range("X4").Value = -1.16
$ws.Range("Y4").Value = -5.35
$ws.Range("Z4").Value = -1.09
$ws.Range("AA4").Value = 292.53
$ws.Range("AB4").Value = 38.05
$ws.Range("AC4").Value = -21
$ws.Range("AD4").Value = -87.77
$ws.Range("AE4").Value = 430
$ws.Range("AF4").Value = 4.31
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 856473009

# Row 5
$ws.Range("D5").Value = 15702
$ws.Range("E5").Value = 496
$ws.Range("F5").Value = 496
$ws.Range("G5").Value = -90
$ws.Range("H5").Value = -455
$ws.Range("I5").Value = -402
$ws.Range("J5").Value = -54
$ws.Range("K5").Value = 13079
$ws.Range("L5").Value = 9731
$ws.Range("M5").Value = 3348
$ws.Range("N5").Value = 3210
$ws.Range("O5").Value = 138
$ws.Range("P5").Value = 4282
$ws.Range("Q5").Value = 319
$ws.Range("R5").Value = -402
$ws.Range("S5").Value = -370
$ws.Range("T5").Value = 214
$ws.Range("U5").Value = 105
$ws.Range("V5").Value = 6433
$ws.Range("W5").Value = 3.16
$ws.Range("X5").Value = -2.9
$ws.Range("Y5").Value = -11.73
$ws.Range("Z5").Value = -3.41
$ws.Range("AA5").Value = 290.67
$ws.Range("AB5").Value = 27.47
$ws.Range("AC5").Value = -47
$ws.Range("AD5").Value = -23.99
$ws.Range("AE5").Value = 379
$ws.Range("AF5").Value = 2.97
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 856473009

# Row 6
$ws.Range("D6").Value = 16488
$ws.Range("E6").Value = 495
$ws.Range("F6").Value = 495
$ws.Range("G6").Value = -14
$ws.Range("H6").Value = -59
$ws.Range("I6").Value = -30
$ws.Range("K6").Value = 11729
$ws.Range("L6").Value = 8510
$ws.Range("M6").Value = 3219
$ws.Range("N6").Value = 3053
$ws.Range("P6").Value = 4282
$ws.Range("Q6").Value = 144
$ws.Range("R6").Value = 582
$ws.Range("S6").Value = -699
$ws.Range("T6").Value = 107
$ws.Range("U6").Value = 37
$ws.Range("V6").Value = 5966
$ws.Range("W6").Value = 3
$ws.Range("X6").Value = -0.36
$ws.Range("Y6").Value = -0.94
$ws.Range("Z6").Value = -0.47
$ws.Range("AA6").Value = 264.36
$ws.Range("AB6").Value = 49.4
$ws.Range("AC6").Value = -3
$ws.Range("AD6").Value = -294.43
$ws.Range("AE6").Value = 361
$ws.Range("AF6").Value = 2.81
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 856473009
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7: clear all data columns, keep A:C
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all data columns, keep A:C
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all data columns, keep A:C
$ws.Range("D9:AJ9").ClearContents()
